$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 148.23077
$ws.Range("I33").Value = 127.25
$ws.Range("K33").Value = 127.25
$ws.Range("M33").Value = 101.75

$ws.Range("H80").Value = 1082.375
$ws.Range("I80").Value = 738.4286
$ws.Range("J80").Value = 1349.8889
$ws.Range("K80").Value = 2215.2858
$ws.Range("L80").Value = 4049.6667
$ws.Range("M80").Value = -1217.2858
$ws.Range("N80").Value = -6045.6667

$ws.Range("H83").Value = 1082.375
$ws.Range("I83").Value = 738.4286
$ws.Range("J83").Value = 1349.8889
$ws.Range("K83").Value = 6645.8574
$ws.Range("L83").Value = 12149.0001
$ws.Range("M83").Value = -1653.8574
$ws.Range("N83").Value = -22133.0001

$ws.Range("H137").Value = 2964.8667
$ws.Range("I137").Value = 1909.6666
$ws.Range("K137").Value = 5728.9998
$ws.Range("M137").Value = -3178.9998

$ws.Range("H138").Value = 3803.5
$ws.Range("I138").Value = 3993.5
$ws.Range("J138").Value = 3613.5
$ws.Range("K138").Value = 11980.5
$ws.Range("L138").Value = 10840.5
$ws.Range("M138").Value = -6840.5
$ws.Range("N138").Value = -21120.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3939.3096
$ws.Range("I32").Value = 3939.3096
$ws.Range("K32").Value = 3939.3096
$ws.Range("M32").Value = -3652.3096

$ws.Range("H61").Value = 8500
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 8500
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 8500
$ws.Range("M61").ClearContents() | Out-Null
$ws.Range("N61").Value = -8924

$ws.Range("H74").Value = 4274.6665
$ws.Range("I74").Value = 1412.5
$ws.Range("K74").Value = 1412.5
$ws.Range("M74").Value = -538.5

$ws.Range("H77").Value = 4274.6665
$ws.Range("I77").Value = 1412.5
$ws.Range("K77").Value = 7062.5
$ws.Range("M77").Value = -2694.5

$ws.Range("H80").Value = 50088
$ws.Range("J80").Value = 50088
$ws.Range("L80").Value = 50088
$ws.Range("N80").Value = -52084

$ws.Range("H83").Value = 50088
$ws.Range("J83").Value = 50088
$ws.Range("L83").Value = 150264
$ws.Range("N83").Value = -160248

$ws.Range("H122").Value = 1363.3572
$ws.Range("I122").Value = 1363.3572
$ws.Range("K122").Value = 4090.0716
$ws.Range("M122").Value = -1640.0716

$ws.Range("H132").Value = 3951.5
$ws.Range("J132").Value = 4989.875
$ws.Range("L132").Value = 14969.625
$ws.Range("N132").Value = -20029.625

$ws.Range("H136").Value = 8500
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 8500
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 25500
$ws.Range("M136").ClearContents() | Out-Null
$ws.Range("N136").Value = -30600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 456.4
$ws.Range("I64").Value = 303
$ws.Range("K64").Value = 303
$ws.Range("M64").Value = -78

$ws.Range("H67").Value = 456.4
$ws.Range("I67").Value = 303
$ws.Range("K67").Value = 303
$ws.Range("M67").Value = 477

$ws.Range("H86").Value = 6222
$ws.Range("I86").Value = 5382.4443
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 5382.4443
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -4259.4443
$ws.Range("N86").Value = -12246

$ws.Range("H89").Value = 6222
$ws.Range("I89").Value = 5382.4443
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 26912.2215
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -21296.2215
$ws.Range("N89").Value = -61232

$ws.Range("H105").Value = 2173.375
$ws.Range("I105").Value = 2126.7144
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 2126.7144
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -379.7143999999998
$ws.Range("N105").Value = -5994

$ws.Range("H132").Value = 150000
$ws.Range("J132").Value = 150000
$ws.Range("L132").Value = 150000
$ws.Range("N132").Value = -160120

$ws.Range("H134").Value = 6831.857
$ws.Range("I134").Value = 6303.8335
$ws.Range("K134").Value = 18911.5005
$ws.Range("M134").Value = -16376.5005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5983.9
$ws.Range("J58").Value = 11236.857
$ws.Range("L58").Value = 11236.857
$ws.Range("N58").Value = -11642.857

$ws.Range("H134").Value = 2395.125
$ws.Range("I134").Value = 2161.1304
$ws.Range("K134").Value = 6483.3912
$ws.Range("M134").Value = -3948.3912

$ws.Range("H136").Value = 5983.9
$ws.Range("J136").Value = 11236.857
$ws.Range("L136").Value = 33710.571
$ws.Range("N136").Value = -38810.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 1900000
$ws.Range("I128").Value = 1900000
$ws.Range("K128").Value = 5700000
$ws.Range("M128").Value = -5695020

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 122682.11
$ws.Range("I132").Value = 154305.72
$ws.Range("K132").Value = 462917.16
$ws.Range("M132").Value = -460387.16

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2337.9167
$ws.Range("I93").Value = 2285.5
$ws.Range("J93").Value = 2600
$ws.Range("K93").Value = 2285.5
$ws.Range("L93").Value = 2600
$ws.Range("M93").Value = -1037.5
$ws.Range("N93").Value = -5096

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 77000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 77000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 77000
$ws.Range("M51").ClearContents() | Out-Null
$ws.Range("N51").Value = -78020

$ws.Range("H81").Value = 799.6667
$ws.Range("I81").Value = 799.6667
$ws.Range("K81").Value = 1599.3334
$ws.Range("M81").Value = -538.3334

$ws.Range("H84").Value = 799.6667
$ws.Range("I84").Value = 799.6667
$ws.Range("K84").Value = 7996.666999999999
$ws.Range("M84").Value = -2692.666999999999

$ws.Range("H100").Value = 434.73334
$ws.Range("I100").Value = 346.375
$ws.Range("J100").Value = 535.7143
$ws.Range("K100").Value = 692.75
$ws.Range("L100").Value = 1071.4286
$ws.Range("M100").Value = -151.75
$ws.Range("N100").Value = -2153.4286

$ws.Range("H136").Value = 3874.8948
$ws.Range("I136").Value = 2973.1538
$ws.Range("J136").Value = 5828.6665
$ws.Range("K136").Value = 8919.4614
$ws.Range("L136").Value = 17485.9995
$ws.Range("M136").Value = -6369.4614
$ws.Range("N136").Value = -22585.9995
